# Applies the cryptos.xlsx price/volume/coin updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be written as TEXT (matches the
# source file, where every cell in B:E is stored as a text/inline string -
# e.g. "1.02" must stay the literal text "1.02", not become the number 1.02).
# Writing straight to the target cell with a "@" text format would leave a
# permanent style-index change on that cell, so instead we stage the text in
# an unused scratch cell and bring it over with Paste Special (values only),
# which copies the text payload but not the formatting.
function Set-TextValue($cellRange, $value) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($cellRange).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "51.355.22"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.075.49"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "394.01"
$ws.Range("E5").Value = "  +2.48%  "
Set-TextValue "D6" "102.51"
$ws.Range("E6").Value = "  -0.39%  "
Set-TextValue "D7" "0.534"
$ws.Range("E7").Value = "  -1.73%  "
Set-TextValue "D9" "0.587"
$ws.Range("E9").Value = "  -0.61%  "
Set-TextValue "D10" "37.45"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +0.78%  "
Set-TextValue "D12" "0.0853"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "3.553.19"
$ws.Range("E13").Value = "  +0.99%  "
Set-TextValue "D14" "18.73"
$ws.Range("E14").Value = "  -0.09%  "
Set-TextValue "D15" "7.70"
$ws.Range("E15").Value = "  -0.85%  "
Set-TextValue "D16" "1.02"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").Value = "3.066.94"
$ws.Range("E17").Value = "  +0.93%  "
Set-TextValue "D18" "10.56"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "51.346.51"
$ws.Range("E19").Value = "  -0.63%  "
Set-TextValue "D20" "3.16"
$ws.Range("E20").Value = "  +1.86%  "
Set-TextValue "D21" "12.34"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  -0.44%  "
Set-TextValue "D23" "70.18"
$ws.Range("E23").Value = "  +0.23%  "
Set-TextValue "D24" "264.88"
$ws.Range("E24").Value = "  -0.90%  "
Set-TextValue "D25" "3.20"
$ws.Range("E25").Value = "  +0.87%  "
Set-TextValue "D26" "7.89"
$ws.Range("E26").Value = "  -6.18%  "
Set-TextValue "D27" "26.98"
$ws.Range("E27").Value = "  +2.14%  "
Set-TextValue "D28" "7.20"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("E29").Value = "  +0.07%  "
Set-TextValue "D30" "0.165"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("E31").Value = "  -2.41%  "
Set-TextValue "D32" "10.71"
$ws.Range("E32").Value = "  +3.96%  "
Set-TextValue "D33" "0.0493"
$ws.Range("E33").Value = "  +10.82%  "
Set-TextValue "D34" "36.27"
$ws.Range("E34").Value = "  +6.32%  "
$ws.Range("E35").Value = "  +0.15%  "
Set-TextValue "D36" "49.92"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -0.11%  "
Set-TextValue "D38" "3.33"
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "4.02"
$ws.Range("E39").Value = "  +9.32%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D40" "0.288"
$ws.Range("E40").Value = "  +0.65%  "
Set-TextValue "D41" "128.89"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  -1.28%  "
Set-TextValue "D43" "16.65"
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D44" "0.115"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D45" "2.53"
$ws.Range("E45").Value = "  +0.02%  "
Set-TextValue "D46" "21.72"
$ws.Range("E46").Value = "  -0.08%  "
Set-TextValue "D47" "2.52"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "2.071.88"
$ws.Range("E49").Value = "  +1.84%  "
Set-TextValue "D50" "0.0498"
$ws.Range("E50").Value = "  +27.18%  "
Set-TextValue "D51" "0.905"
$ws.Range("E51").Value = "  +10.50%  "

# Clean up the scratch cell and clipboard marching ants.
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

